$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9333678483963013
$ws.Range("B1").Value = 0.8640801310539246
$ws.Range("D1").Value = 1.568647027015686
$ws.Range("E1").Value = 0.9496013522148132
